$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 69-74 (Aged Care outbreak / resident / staff case metrics) are unchanged.
# Rows 75-78 previously held the "# / % Aged Care Molnupiravir/Paxlovid Prescriptions (Daily)"
# metrics (6 rows, 75-80). They are replaced by 4 "Residential Aged Care" branded
# metrics (no more "per 1M" variants), and the "# Monthly PBS Scripts" trio shifts
# up from rows 81-83 to rows 79-81. The sheet therefore shrinks from 83 to 81 rows.

$ws.Range("C75").Value = "# Residential Aged Care Molnupiravir Prescriptions (Daily)"

$ws.Range("C76").Value = "% Residential Aged Care Molnupiravir Prescriptions (Daily) per Case"
$ws.Range("D76").Value = 790

$ws.Range("C77").Value = "# Residential Aged Care Paxlovid Prescriptions (Daily)"
$ws.Range("D77").Value = 800

$ws.Range("C78").Value = "% Residential Aged Care Paxlovid Prescriptions (Daily) per Case"
$ws.Range("D78").Value = 820

$ws.Range("C79").Value = "# Monthly PBS Scripts"
$ws.Range("D79").Value = 830
$ws.Range("C79").Style = "Normal"

$ws.Range("C80").Value = "# Monthly PBS Scripts per 1M"
$ws.Range("D80").Value = 840
$ws.Range("C80").Style = "Normal"

$ws.Range("C81").Value = "% Monthly PBS Scripts Change"
$ws.Range("D81").Value = 850

# The old rows 82 and 83 (formerly "# Monthly PBS Scripts per 1M" / "% Monthly PBS
# Scripts Change") are no longer needed now that the table is 2 rows shorter.
$ws.Rows("82:83").Delete()

# Column C now holds longer "Residential Aged Care ..." labels; re-fit the column.
$ws.Columns("C:C").AutoFit()

$ws.Range("C78").Select()
